# Sincronizacao de dados: novo orcamento (Douglas Hiromitsu) inserido no topo
# da lista de "quotations" (nova linha 18), com seus dois itens de pedido
# (visita tecnica + km rodado) inseridos no topo do bloco correspondente em
# "items" (novas linhas 66 e 67). Todas as linhas subsequentes das duas
# planilhas se deslocam para baixo, preservando seu conteudo original.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "quotations": insert a new row 18 with the new quotation record.
# ---------------------------------------------------------------------
$qs = $wb.Worksheets.Item("quotations")
$qs.Rows.Item(18).Insert()

$qRow = $qs.Range("A18:U18")
$qRow.NumberFormat = "@"

# id of the newly inserted quotation, reused below for the "items" rows
$newQuotationId = "NmIyYjA3ZmYtY2I3MC00NTBiLWJjMDMtNTM4ZGEyNmI1ZmIwOjU3MDE2"

$qs.Cells.Item(18, 1).Value  = $newQuotationId
$qs.Cells.Item(18, 2).Value  = "HJCRTZS--G"
$qs.Cells.Item(18, 3).Value  = "31415680817 - DOUGLAS HIROMITSU TFK 2Kg"
$qs.Cells.Item(18, 4).Value  = "Parafuso (Pino) que sustenta raspador quebrou, iremos fazer uma visita técnica para descobrir se há mais defeitos. 11 96619-3370"
$qs.Cells.Item(18, 5).Value  = ""
$qs.Cells.Item(18, 6).Value  = $false
$qs.Cells.Item(18, 7).Value  = "484"
$qs.Cells.Item(18, 8).Value  = "484"
$qs.Cells.Item(18, 9).Value  = "Aprovada"
$qs.Cells.Item(18, 10).Value = "2025-12-08T12:21:25.404Z"
$qs.Cells.Item(18, 11).Value = ""
$qs.Cells.Item(18, 12).Value = "douglas simao"
$qs.Cells.Item(18, 13).Value = "Kaue Teixeira Caldeira Venâncio"
$qs.Cells.Item(18, 14).Value = "accounts/57016/quotations/6b2b07ff-cb70-450b-bc03-538da26b5fb0/signatures/82045e1d-cc6f-45af-a3f3-2c6a0906455d.png"
$qs.Cells.Item(18, 15).Value = "2025-12-01T12:28:02.404Z"
$qs.Cells.Item(18, 16).Value = ""
$qs.Cells.Item(18, 17).Value = "percentage"
$qs.Cells.Item(18, 18).Value = "0"
$qs.Cells.Item(18, 19).Value = "0"
$qs.Cells.Item(18, 20).Value = "NDUwNzM2Mzo1NzAxNg=="
$qs.Cells.Item(18, 21).Value = "approved"

# ---------------------------------------------------------------------
# Sheet "items": insert two new rows (66 and 67) with the two new items
# belonging to the quotation inserted above.
# ---------------------------------------------------------------------
$is_ = $wb.Worksheets.Item("items")
$is_.Range("A66:J67").Insert()

# Row 66 - "67Km de distancia"
$is_.Range("A66").NumberFormat = "@"
$is_.Range("D66").NumberFormat = "@"
$is_.Range("F66").NumberFormat = "@"
$is_.Range("G66").NumberFormat = "@"
$is_.Range("I66").NumberFormat = "@"
$is_.Range("J66").NumberFormat = "@"

$is_.Cells.Item(66, 1).Value = "NzVmMDU4ODAtZTJiMi00MmFlLTgzZmItYTJlYmYwY2EyNTgxOjU3MDE2"
$is_.Cells.Item(66, 2).Value = 67
$is_.Cells.Item(66, 3).Value = 134
$is_.Cells.Item(66, 4).Value = "67Km de distancia"
$is_.Cells.Item(66, 5).Value = 3
$is_.Cells.Item(66, 6).Value = $newQuotationId
$is_.Cells.Item(66, 7).Value = "Zjc3ODdmZmQtNzZiNy00ZjNmLThmNjQtNjdjOGIyOGYxYzUwOjU3MDE2"
$is_.Cells.Item(66, 8).Value = 2
$is_.Cells.Item(66, 9).Value = "service"
$is_.Cells.Item(66, 10).Value = $newQuotationId

# Row 67 - "Visita técnica"
$is_.Range("A67").NumberFormat = "@"
$is_.Range("D67").NumberFormat = "@"
$is_.Range("F67").NumberFormat = "@"
$is_.Range("G67").NumberFormat = "@"
$is_.Range("I67").NumberFormat = "@"
$is_.Range("J67").NumberFormat = "@"

$is_.Cells.Item(67, 1).Value = "ZDQ0NDkzNjgtNjE2MS00YjRkLWEyNWUtOTdkOGFiMjFiMjVjOjU3MDE2"
$is_.Cells.Item(67, 2).Value = 1
$is_.Cells.Item(67, 3).Value = 350
$is_.Cells.Item(67, 4).Value = "Visita técnica"
$is_.Cells.Item(67, 5).Value = 3
$is_.Cells.Item(67, 6).Value = $newQuotationId
$is_.Cells.Item(67, 7).Value = "NWVmNmQ0MDEtNzBmMy00Yzg3LWFlZDAtYzJiYTM1MTc4OWNlOjU3MDE2"
$is_.Cells.Item(67, 8).Value = 350
$is_.Cells.Item(67, 9).Value = "service"
$is_.Cells.Item(67, 10).Value = $newQuotationId
